$wb = $excel.ActiveWorkbook

# --- Update the "Conversión del día" message on sheet "Hoja1" ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 1.82 = 6592.82 pesos`n✅ 6592.82 pesos = 1.81 = 940.16 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$wsHoja1.Range("A1").Value = $newText

# --- Update the rate figures on sheet "tasas" ---
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("N10").Value = 550
$wsTasas.Range("O10").Value = 3626.05

$wsTasas.Range("N12").Value = 3639.99
$wsTasas.Range("O12").Value = 519.078
